$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 73: 2014-10-05, 22:59 -> 23:35, 5 min interruption, category "Coding"
$ws.Range("A73").Value = 41917
$ws.Range("B73").Value = 0.95763888888888893
$ws.Range("C73").Value = 0.98263888888888884
$ws.Range("D73").Value = 5
$ws.Range("F73").Value = "Coding"
$ws.Range("E73").Formula = '=IF(AND(NOT(ISBLANK(B73)),NOT(ISBLANK(C73))),(C73-B73)*24-D73/60,"")'

# Row 74: 2014-10-06, 13:20 -> 14:12, 15 min interruption, category "Coding"
$ws.Range("A74").Value = 41918
$ws.Range("B74").Value = 0.55555555555555558
$ws.Range("C74").Value = 0.59166666666666667
$ws.Range("D74").Value = 15
$ws.Range("F74").Value = "Coding"
$ws.Range("E74").Formula = '=IF(AND(NOT(ISBLANK(B74)),NOT(ISBLANK(C74))),(C74-B74)*24-D74/60,"")'

# Move the active selection to A75, matching where entry would continue next
$ws.Range("A75").Select() | Out-Null
